$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# -----------------------------------------------------------------
# 1) Merge the two runs ", Dominik, " + "if24b161@technikum-wien.at"
#    into a single run (same formatting, so a same-text Find/Replace
#    scoped to that paragraph makes Word re-flow it into one run).
# -----------------------------------------------------------------
$pBiesaga = $d.Paragraphs.Item(10)
$rBiesaga = $pBiesaga.Range
$rBiesaga.Find.Execute(", Dominik, if24b161@technikum-wien.at", $true, $false, $false, $false, $false, $true, 1, $false, ", Dominik, if24b161@technikum-wien.at", 2) | Out-Null

# -----------------------------------------------------------------
# 2) Replace the empty paragraph after "Problem description" with two
#    new German paragraphs describing the problem.
# -----------------------------------------------------------------
$pEmpty = $d.Paragraphs.Item(16)
$rEmpty = $pEmpty.Range
$rEmpty.Collapse(1)
$xmlProblem = "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"de-AT`"/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:lang w:val=`"de-AT`"/></w:rPr><w:t>Da</w:t></w:r>" + `
  "<w:r><w:rPr><w:lang w:val=`"de-AT`"/></w:rPr><w:t>s Studententeam braucht Abwechslung in den Pausen zwischen den Vorlesungen, und benötigt eine neue Interpretation von einem klassischen Spiel.</w:t></w:r>" + `
  "</w:p>" + `
  "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"de-AT`"/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:lang w:val=`"de-AT`"/></w:rPr><w:t xml:space=`"preserve`">Unser Team möchte neue Erfahrungen im Bereich Spielentwicklung und Game </w:t></w:r>" + `
  "<w:proofErr w:type=`"spellStart`"/>" + `
  "<w:r><w:rPr><w:lang w:val=`"de-AT`"/></w:rPr><w:t>Engines</w:t></w:r>" + `
  "<w:proofErr w:type=`"spellEnd`"/>" + `
  "<w:r><w:rPr><w:lang w:val=`"de-AT`"/></w:rPr><w:t xml:space=`"preserve`"> sammeln.</w:t></w:r>" + `
  "</w:p>"
$rEmpty.InsertXML($xmlProblem)

# -----------------------------------------------------------------
# 3) Split "Solution description" into "Solution " + "description"
#    runs, wrapping "description" in a spellcheck proofErr pair.
#    (One extra paragraph was inserted above, so this heading is now
#    paragraph 18.)
# -----------------------------------------------------------------
$pSolution = $d.Paragraphs.Item(18)
$rSolution = $pSolution.Range
$xmlSolution = "<w:p $wns><w:pPr><w:pStyle w:val=`"Heading1`"/>" + `
  "<w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"D9D9D9`" w:themeFill=`"background1`" w:themeFillShade=`"D9`"/>" + `
  "<w:tabs><w:tab w:val=`"left`" w:pos=`"2775`"/></w:tabs>" + `
  "<w:ind w:firstLine=`"142`"/><w:jc w:val=`"left`"/>" + `
  "<w:rPr><w:rFonts w:asciiTheme=`"minorHAnsi`" w:hAnsiTheme=`"minorHAnsi`" w:cstheme=`"minorHAnsi`"/><w:b w:val=`"0`"/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/><w:lang w:val=`"de-AT`"/></w:rPr>" + `
  "</w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:asciiTheme=`"minorHAnsi`" w:hAnsiTheme=`"minorHAnsi`" w:cstheme=`"minorHAnsi`"/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/><w:lang w:val=`"de-AT`"/></w:rPr><w:t xml:space=`"preserve`">Solution </w:t></w:r>" + `
  "<w:proofErr w:type=`"spellStart`"/>" + `
  "<w:r><w:rPr><w:rFonts w:asciiTheme=`"minorHAnsi`" w:hAnsiTheme=`"minorHAnsi`" w:cstheme=`"minorHAnsi`"/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/><w:lang w:val=`"de-AT`"/></w:rPr><w:t>description</w:t></w:r>" + `
  "<w:proofErr w:type=`"spellEnd`"/>" + `
  "</w:p>"
$rSolution.InsertXML($xmlSolution)

# -----------------------------------------------------------------
# 4) Merge the two runs "Godot game engin" + "e" into a single run
#    (same formatting). Scope the Find to the "Technical environment"
#    paragraph that reads exactly "Godot game engine" so the other
#    "Godot ... engine" mention in the solution text is untouched.
# -----------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Godot game engine`r") {
        $rGodot = $p.Range
        $rGodot.Find.Execute("Godot game engine", $true, $false, $false, $false, $false, $true, 1, $false, "Godot game engine", 2) | Out-Null
        break
    }
}
